# تعديل تلقائي في شيت Card3 by admin at 2025-12-06 18:34:12
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card3")

# Update the "card" column (A) from 2 -> 3 for rows 3-7 and 9-13
# (rows 2 and 8 already contain 3 and are left untouched)
# Leading apostrophe keeps the value stored as text "3", matching the
# original inline-string cell type rather than converting it to a number.
$rows = @(3, 4, 5, 6, 7, 9, 10, 11, 12, 13)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 1).Value = "'3"
}
